$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 640; existing rows 640-720 shift down to 641-721.
$ws.Rows(640).Insert()

# Populate the newly inserted row 640 with the new data record.
$ws.Cells.Item(640, 1).Value = 5
$ws.Cells.Item(640, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(640, 3).Value = "Maule"
$ws.Cells.Item(640, 4).Value = 44918
$ws.Cells.Item(640, 5).Value = 7
$ws.Cells.Item(640, 6).Value = 100112002
$ws.Cells.Item(640, 7).Value = "Pimiento"
$ws.Cells.Item(640, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(640, 9).Value = "Primera"
$ws.Cells.Item(640, 10).Value = 400
$ws.Cells.Item(640, 11).Value = 8000
$ws.Cells.Item(640, 12).Value = 8000
$ws.Cells.Item(640, 13).Value = 8000
$ws.Cells.Item(640, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(640, 15).Value = "Región del Maule"
$ws.Cells.Item(640, 16).Value = 444
$ws.Cells.Item(640, 17).Value = 18
$ws.Cells.Item(640, 18).Value = "Hortaliza"
